$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Contestants" sheet: rows 2-4 cyclically rotate their contents.
#   row2 -> row3, row3 -> row4, row4 -> row2
# (This reflects the underlying seats/contestants being reshuffled.)
# ---------------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

# Row 2 (becomes old row 4's data: Kathleen Reynolds)
$contestants.Range("A2").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$contestants.Range("B2").Value = "Kathleen Reynolds"
$contestants.Range("C2").Value = 33
$contestants.Range("E2").Value = "kathleenmonicareynolds@gmail.com"
$contestants.Range("G2").Value = "Footscray"
$contestants.Range("J2").Value = "Peter Adamidis, Felicity Parker-Hill"
$contestants.Range("L2").Value = "N"
$contestants.Range("M2").Value = "N/A"

# Row 3 (becomes old row 2's data: Felicity Parker-Hill)
$contestants.Range("A3").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$contestants.Range("B3").Value = "Felicity Parker-Hill"
$contestants.Range("C3").Value = 27
$contestants.Range("E3").Value = "felicity.parkerhill@endemolshine.com.au"
$contestants.Range("G3").Value = "Melbourne"
$contestants.Range("J3").Value = "Peter Adamidis, Kathleen Reynolds"
$contestants.Range("L3").Value = "N"
$contestants.Range("M3").Value = "N/A"

# Row 4 (becomes old row 3's data: Peter Adamidis)
$contestants.Range("A4").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$contestants.Range("B4").Value = "Peter Adamidis"
$contestants.Range("C4").Value = 34
$contestants.Range("E4").Value = "peter.adamidis@gmail.com"
$contestants.Range("G4").Value = ""
$contestants.Range("J4").Value = "Kathleen Reynolds, Felicity Parker-Hill"
$contestants.Range("L4").Value = "Y"
$contestants.Range("M4").Value = "Broken Leg"

# ---------------------------------------------------------------------------
# "Seat Assignments" sheet: each row gets a fresh assignment ID and its
# ContestantID / Seat updated to reflect the new seating. Row 4's stray
# empty Notes cell (H4) is also cleared out entirely.
# ---------------------------------------------------------------------------
$seats = $wb.Worksheets.Item("Seat Assignments")

$seats.Range("A2").Value = "20946cbf-8f88-4629-a6d2-d8510554cb9a"
$seats.Range("C2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$seats.Range("E2").Value = "C1"

$seats.Range("A3").Value = "c1d29895-e1ab-4926-8d3b-29a9df7058de"
$seats.Range("E3").Value = "C3"

$seats.Range("A4").Value = "f8463207-3f82-429c-bf09-986ae9a6cc97"
$seats.Range("C4").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$seats.Range("E4").Value = "D3"
$seats.Range("H4").ClearContents()
